$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.692.22"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "3.072.33"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'551.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").Value = "'141.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.02%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.067.13"
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("D9").Value = "'0.503"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "'6.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.14%  "
$ws.Range("D11").Value = "'0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").Value = "'35.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "3.566.99"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").Value = "63.621.82"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("D17").Value = "3.072.08"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "'6.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").Value = "'487.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.71%  "
$ws.Range("D21").Value = "'13.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("D22").Value = "'0.679"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'7.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.68%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").Value = "'12.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.12%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").Value = "'7.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").Value = "'2.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.32%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "'26.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +6.92%  "
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("D35").Value = "'55.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("D37").Value = "'469.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("E38").Value = "  +4.48%  "
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("D40").Value = "3.072.78"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").Value = "'8.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").Value = "'2.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.51%  "
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("D45").Value = "'0.258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.67%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").Value = "0.0₃0518"
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("D50").Value = "'117.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").Value = "'2.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.40%  "
